# Updated cryptos list (Price / Volume(1h) columns) per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> @(newPriceTextOrNull, newVolumeText)
$rowData = @{
    2 = @("27.908.87", "  -0.03%  ")
    3 = @("1.632.20", "  -0.79%  ")
    4 = @("0.998", "  -0.26%  ")
    5 = @("211.68", "  -0.88%  ")
    6 = @("0.523", "  -0.76%  ")
    7 = @("0.997", "  -0.28%  ")
    8 = @($null, "  -0.59%  ")
    9 = @($null, "  -2.44%  ")
    10 = @($null, "  -0.18%  ")
    11 = @($null, "  +1.01%  ")
    12 = @("1.860.33", "  -1.00%  ")
    13 = @("1.626.37", "  -1.19%  ")
    14 = @("4.05", "  -0.62%  ")
    15 = @($null, "  -0.14%  ")
    16 = @("65.36", "  -0.32%  ")
    17 = @("27.894.78", "  -0.09%  ")
    18 = @("230.68", "  -0.25%  ")
    19 = @($null, "  -0.15%  ")
    20 = @($null, "  -1.73%  ")
    21 = @("0.997", "  -0.37%  ")
    22 = @($null, "  -0.63%  ")
    23 = @("10.35", "  -2.86%  ")
    24 = @($null, "  -3.82%  ")
    25 = @($null, "  +1.48%  ")
    26 = @("6.95", "  +0.51%  ")
    27 = @($null, "  -0.92%  ")
    28 = @("15.62", "  -0.73%  ")
    29 = @("0.998", "  -0.44%  ")
    31 = @($null, "  -0.76%  ")
    32 = @("3.42", "  +2.49%  ")
    33 = @("1.408.28", "  -2.39%  ")
    34 = @("3.08", "  +0.32%  ")
    35 = @($null, "  +0.09%  ")
    36 = @("1.02", "  +9.03%  ")
    37 = @($null, "  +1.11%  ")
    38 = @($null, "  +0.56%  ")
    39 = @("0.560", "  +0.24%  ")
    40 = @("0.868", "  -2.39%  ")
    41 = @($null, "  +0.00%  ")
    42 = @("0.997", "  -0.35%  ")
    43 = @("66.63", "  -3.59%  ")
    44 = @($null, "  +1.73%  ")
    45 = @("1.83", "  +0.09%  ")
    46 = @("2.20", "  -0.81%  ")
    47 = @("1.772.20", "  -0.92%  ")
    48 = @("87.88", "  -1.31%  ")
    49 = @($null, "  -1.11%  ")
    50 = @("0.0998", "  -1.05%  ")
    51 = @("0.0507", "  -0.21%  ")
}

# Cell used as a known "clean" (no explicit style) reference so that forcing
# column D values to be stored as text does not leave a residual style on the cell.
$cleanStyle = $ws.Range("D30").Style

foreach ($r in $rowData.Keys) {
    $pair = $rowData[$r]
    $priceText = $pair[0]
    $volumeText = $pair[1]
    if ($priceText -ne $null) {
        $priceCell = $ws.Range("D$r")
        $priceCell.NumberFormat = "@"
        $priceCell.Value = $priceText
        $priceCell.Style = $cleanStyle
    }
    $ws.Range("E$r").Value = $volumeText
}
